# Actualiza la tabla de resultados (Practica_3) con los nuevos datos del informe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Columna B (findMaxElement) y C (findMaxElementBySorting) ---
# Fila 3 (Tamaño 100)
$ws.Range("C3").Value = 0

# Fila 4 (Tamaño 1000)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1

# Fila 5 (Tamaño 1000000)
$ws.Range("B5").Value = 60
$ws.Range("C5").Value = 1423

# Fila 6 (Tamaño 2000000)
$ws.Range("B6").Value = 39
$ws.Range("C6").Value = 872

# Fila 7 (Tamaño 2500000)
$ws.Range("B7").Value = 24
$ws.Range("C7").Value = 926

# Fila 8 (Tamaño 5000000)
$ws.Range("B8").Value = 203
$ws.Range("C8").Value = 5324

# Fila 9 (Tamaño 6000000)
$ws.Range("B9").Value = 83
$ws.Range("C9").Value = 2626

# Fila 10 (Tamaño 6500000)
$ws.Range("C10").Value = 2839

# Fila 11 (Tamaño 7000000)
$ws.Range("B11").Value = 103
$ws.Range("C11").Value = 3243

# Fila 12 (Tamaño 75000000)
$ws.Range("B12").Value = 1451
$ws.Range("C12").Value = 45621

# --- Formato: A6, A9, A10 y A11 pierden el borde "con numFmt" y pasan a
#     compartir el estilo normal de borde (igual que A2 / A7) ---
$ws.Range("A6").Borders.LineStyle = 1
$ws.Range("A6").Borders.Weight = 2
$ws.Range("A6").Font.Name = "Aptos"
$ws.Range("A6").Font.Size = 11
$ws.Range("A6").HorizontalAlignment = -4130
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeRight=10 (sin xlEdgeBottom=9: este estilo
# deja el borde inferior abierto, igual que A7/A8)
foreach ($cell in @("A9", "A10", "A11")) {
    $rng = $ws.Range($cell)
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(10).Weight = 2
    $rng.Font.Name = "Aptos"
    $rng.Font.Size = 11
    $rng.HorizontalAlignment = -4130
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# --- Selección final tal y como queda guardada en el libro ---
$ws.Range("A12:C12").Select()
